# Update Excel files from OneDrive - Tue May 27 05:55:35 UTC 2025
#
# The "Test Data" sheet gets the Test-Start / Test-End dates filled in for
# rows 2-5, together with the Test-Result / Status drop-down picks that go
# with them. The user's cursor ends up on F4 (Test Data) after doing this,
# having scrolled the sheet one column to the right (topLeftCell = B1), and
# on G12 of the "Ref" sheet (left over from picking values out of its lookup
# tables) while "Test Data" remains the active tab.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Test Data")
$wsRef  = $wb.Worksheets.Item("Ref")

$wsData.Activate()

# Row 2 - ENCAPSULANT / GEL TEST / HIUV
$wsData.Range("E2").Value = 45804
$wsData.Range("F2").Value = 45804
$wsData.Range("G2").Value = "Pending"
$wsData.Range("H2").Value = "Pending"

# Row 3 - ENCAPSULANT / SHRINKAGE TEST / HIUV
$wsData.Range("E3").Value = 45802
$wsData.Range("F3").Value = 45802
$wsData.Range("G3").Value = "Pass"
$wsData.Range("H3").Value = "Completed"

# Row 4 - ENCAPSULANT / GEL TEST / SUNBEZ
$wsData.Range("E4").Value = 45804
$wsData.Range("F4").Value = 45804
$wsData.Range("G4").Value = "Pending"
$wsData.Range("H4").Value = "Pending"

# Row 5 - ENCAPSULANT / SHRINKAGE TEST / SUNBEZ
$wsData.Range("E5").Value = 45802
$wsData.Range("F5").Value = 45802
$wsData.Range("G5").Value = "Pass"
$wsData.Range("H5").Value = "Completed"

# Leftover selection on the "Ref" sheet from picking the drop-down values
# above (sheet stays inactive; only its remembered selection changes).
$wsRef.Activate()
$wsRef.Range("G12").Select()

# Back to "Test Data" as the active tab, scrolled right one column with the
# cursor resting on F4.
$wsData.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$wsData.Range("F4").Select()
